# Implemented VS Studio Integration with Unity, Updated Story
#
# The backlog items that used to live in rows 6-12 move down two rows
# (to rows 8-14); a couple of new items fill in the vacated space at the
# bottom (rows 13-15, previously blank placeholders) and rows 6-7 are
# cleared out (row 6 keeps its Priority flag, row 7 becomes fully blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: clear the Epic/Description/Priority/Vertics/SP content, keep B6=1 ---
$ws.Range("C6:F6").Clear()

# --- Row 7: clear everything (including the Priority flag in B7) ---
$ws.Range("B7:F7").Clear()
$ws.Rows.Item(7).AutoFit()

# --- Row 8: "Simple Level Structure of the first Level" item ---
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "Simple Level Structure of the first Level"
$ws.Range("D8").Value = "!!!"
$ws.Range("E8").Value = "As a player I want to move around in the second level. To progress in the story."
$ws.Range("F8").Value = 5

# --- Row 9: "Design new Enemies" item (note: gets the taller 30pt row) ---
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Design new Enemies"
$ws.Range("D9").Value = "!"
$ws.Range("E9").Value = "As a Player I want to fight against 100  enemies. That means 10 new enemies per level. I also want to see some familiar enemies which I already met in previous games."
$ws.Range("F9").Value = 5
$ws.Rows.Item(9).RowHeight = 30

# --- Row 10: "Convert Level form Paper to Unity" item (loses the 30pt row) ---
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Convert Level form Paper to Unity"
$ws.Range("D10").Value = "!"
$ws.Range("E10").Value = "As a player I want to have a good structured level, where I can move around freely and fight enemies."
$ws.Range("F10").Value = 5
$ws.Rows.Item(10).AutoFit()

# --- Row 11: "Design Dialoges" item ---
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Design Dialoges"
$ws.Range("D11").Value = "!"
$ws.Range("E11").Value = "As a player I want to have interesting and hilarious dialoges between the characters.  "
$ws.Range("F11").Value = 2

# --- Row 12: "Implement Questlog" item (gets the taller 30pt row) ---
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Implement Questlog"
$ws.Range("D12").Value = "!"
$ws.Range("E12").Value = 'As a player I want to have a widget to access all my quests in order to have a good overview of all my open quests. Furthermore I want to get informed if a new quest comes up (Display Box for "You got a new Quest")'
$ws.Range("F12").Value = 5
$ws.Rows.Item(12).RowHeight = 30

# --- Row 13: "Old Lady Quest" item (new content, row was previously blank) ---
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Old Lady Quest"
$ws.Range("D13").Value = "!!"
$ws.Range("E13").Value = 'As a player I want to get  the quest to help the old lady bring her "groceries" to her hut.'
$ws.Range("F13").Value = 5

# --- Row 14: "Hut of the old Lady" item (new content, row was previously blank) ---
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Hut of the old Lady"
$ws.Range("D14").Value = "!!"
$ws.Range("E14").Value = "As a player I want to see the hut of the old Lady. Furthermore I want to move around in it."
$ws.Range("F14").Value = 2

# --- Row 15: just gets the Priority flag, rest stays blank ---
$ws.Range("B15").Value = 1

# --- Update the view: drop the D1 scroll anchor, select C6 instead of F12 ---
[void]$ws.Range("C6").Select()
